# Append two new rows (115 and 116) to Sheet1, continuing the daily series
# found in row 114 (same values in columns B:J, next calendar days in A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 115: copy formatting/values from row 114, then overwrite the date (col A).
$ws.Range("A114:J114").Copy($ws.Range("A115:J115"))
$ws.Range("A115").Value2 = 45671

# Row 116: same treatment, next day again.
$ws.Range("A114:J114").Copy($ws.Range("A116:J116"))
$ws.Range("A116").Value2 = 45672
